# "Generate Report for Handoff" - drop the a9377a48... row (handback not yet
# in sync) and refresh the 425d36a0... row's status/timestamps to reflect a
# fresh handoff, on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-38-17 18:38:41"
$ov.Rows(3).Delete()

# Re-seat the surviving hyperlink (the bulk Delete()/Add() dance is the only
# reliable way this host lets us drop the stale row-3 hyperlink entries) so
# the link target for row 2 stays intact.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b0cab221b62bc804ea79e342b135318120eb2ef9/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File
# | Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-17 18:38:38"
$zh.Range("G2").Value = "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf"
$zh.Range("G2").Value = "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf"
$zh.Rows(3).Delete()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b0cab221b62bc804ea79e342b135318120eb2ef9/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b0cab221b62bc804ea79e342b135318120eb2ef9/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cdfb3a911b506b2a91c38b24982eaf9627d4f3a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/062a87a31b674536ce0fc8e446c7f9214f5af277/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c31d048927ee66f7b42ae1853437e5925acfa1bf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same columns as zh-cn
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-17 18:38:41"
$de.Rows(3).Delete()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b0cab221b62bc804ea79e342b135318120eb2ef9/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b0cab221b62bc804ea79e342b135318120eb2ef9/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c03cb625717c647f67a851e27a11e9f9b48dcd38/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/85c1710781d410278b1a77a0aaeec54b745595b7/e2e/425d36a0-c691-45e2-ab11-0a37ec306df5.md", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b7e88ab44a7bb3136ee1bd055e13730fcec00741/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf", "", "", "425d36a0-c691-45e2-ab11-0a37ec306df5.30a09fc566298713cba6fd8bbaf67821415842ff.de-de.xlf") | Out-Null
